# Apply text corrections described in the commit diff.
# These three shared-string texts are used on both test-case blocks
# (TC1 rows 6-14 and TC2 rows 17-24), so update every matching cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    "Lider de Pessoas esta autenticado no sistema; e, tem permissao para gerenciar Metas de Desempenho" = "Lider de Pessoas esta autenticado no sistema e tem permissao para gerenciar Metas de Desempenho"
    "Lider de Pessoas com uma avaliacao selecionada, clica na opcao 'Editar' modificar a Avaliacao de Desempenho" = "Lider de Pessoas com uma avaliacao selecionada, clica na opcao 'Editar' para modificar a Avaliacao de Desempenho"
    "SYSTEM apresenta o formulario com o campo 'Metas' constando cada Competencia do perfil avaliado" = "SYSTEM apresenta o formulario com o campo 'Metas' contendo cada Competencia do perfil avaliado"
}

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null -and $replacements.ContainsKey($val)) {
            $cell.Value = $replacements[$val]
        }
    }
}
